# Tabelle.xlsx — "Add files via upload" edit
#
# The uploaded workbook revision:
#   * drops two shared strings that are no longer referenced as their own
#     distinct concepts ("Lavora_su" and "Omesso") from the ER-diagram
#     glossary sheet, and introduces a single "-" placeholder string used
#     everywhere "Omesso" used to be used;
#   * the one cell that used to read "Lavora_su" (H29, the "Partecipazioni
#     esterne" domain note for Operatore.CodiceO) now reads "Problema"
#     instead;
#   * the current selection/scroll position moved from N27 to N37 (ish).
#
# Everything else in the grid (all the other shared-string index churn in
# the raw XML) is pure artefact of the string table being rebuilt after the
# two removals — the actual cell text is unchanged there, so there is
# nothing to redo for those cells through the object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells that used to display "Omesso" now display "-" ------------------
$omessoCells = @("F4", "F5", "F15", "F16", "F23", "F30", "N30", "F31", "N31", "F32", "N32", "F33", "N33", "N37")
foreach ($ref in $omessoCells) {
    $ws.Range($ref).Value = "-"
}

# --- The one cell that used to display "Lavora_su" now displays "Problema" -
$ws.Range("H29").Value = "Problema"

# --- Selection / scroll moved down the sheet (was N27) --------------------
$ws.Range("N37").Select()
